$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SQL queries stored in column B (rows 2-7) and C2 all joined tables using
# the generic "id" columns (std.id / prt.id) aliased with dotted "table.id"
# style names. The commit updates every one of these queries to use the real,
# fully-qualified key column names (std.study_id / prt.participant_id, etc.)
# on both sides of each LEFT JOIN.

$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cells) {
    $range = $ws.Range($addr)
    $text = $range.Value2

    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $range.Value = $text
}
